$wb = $excel.ActiveWorkbook

# Select the Input sheet and set H2 to "test"
$ws = $wb.Worksheets.Item("Input")
$ws.Activate()
$ws.Range("H2").Value = "test"
$ws.Range("H2").Select()
